# Applies the OOXML diff:
#  - "Output Parameters: " paragraph gains a new run "Loads maps an enemies"
#  - "Modules Called: ... Map 2.0" paragraph loses the _GoBack bookmark
#  - "Author: " -> "Author: Delmis Spies" (single run)
#  - "Date:  " -> "Date:  12/2/2015" (single run)
#  - "Peer Reviewer" + ": " runs collapse into one "Peer Reviewer: Jan Cajas"
#    run, and the _GoBack bookmark is (re)created at the end of that paragraph

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($findText, $innerXml, $pPr) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
        return
    }
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wNs + '><w:body>' +
           '<w:p ' + $pPr + '>' + $innerXml + '</w:p>' +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

$paraAttrs = 'w:rsidR="009D64CA" w:rsidRDefault="009D64CA" w:rsidP="00A23CB3"'

# 1. Output Parameters: add a new run with "Loads maps an enemies"
Set-ParagraphXml "Output Parameters: " `
    ('<w:r><w:t xml:space="preserve">Output Parameters: </w:t></w:r>' +
     '<w:r><w:t>Loads maps an enemies</w:t></w:r>') `
    $paraAttrs

# 2. Modules Called: strip out the _GoBack bookmark
Set-ParagraphXml "Modules Called: User Enter Section of Map 2.0" `
    ('<w:r><w:t xml:space="preserve">Modules Called: </w:t></w:r>' +
     '<w:r><w:t>User Enter Section of Map 2.0</w:t></w:r>') `
    $paraAttrs

# 3. Author
Set-ParagraphXml "Author: " `
    '<w:r><w:t>Author: Delmis Spies</w:t></w:r>' `
    'w:rsidR="009D64CA" w:rsidRDefault="009D64CA" w:rsidP="00AF07CB"'

# 4. Date
Set-ParagraphXml "Date:  " `
    '<w:r><w:t>Date:  12/2/2015</w:t></w:r>' `
    $paraAttrs

# 5. Peer Reviewer: merge the two runs and move the _GoBack bookmark here
Set-ParagraphXml "Peer Reviewer: " `
    ('<w:r><w:t>Peer Reviewer: Jan Cajas</w:t></w:r>' +
     '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>') `
    $paraAttrs

Write-Host "done"
